# Updated cryptos list on Sat Feb 17 17:28:37 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto table on the active sheet, and swaps the VeChain/RenderToken rows
# (34 and 35) which changed rank order in the source feed.
#
# Price values are plain text in this sheet (note the "1.234.56"-style
# thousands separators, which Excel cannot parse as a number anyway) -
# a leading apostrophe is used for the ones that DO look numeric so Excel
# stores them as text instead of silently coercing them to a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Address,
        [string]$Text
    )
    $ws.Range($Address).Value = "'" + $Text
}

# Row 2 - Bitcoin
Set-TextCell "D2" "51.158.05"
$ws.Range("E2").Value = "  -1.64%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.763.98"
$ws.Range("E3").Value = "  -0.51%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextCell "D5" "353.04"
$ws.Range("E5").Value = "  -1.26%  "

# Row 6 - Solana
Set-TextCell "D6" "108.03"
$ws.Range("E6").Value = "  -1.25%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -2.64%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.01%  "

# Row 9 - Cardano
Set-TextCell "D9" "0.582"
$ws.Range("E9").Value = "  -1.51%  "

# Row 10 - Avalanche
Set-TextCell "D10" "39.43"
$ws.Range("E10").Value = "  -1.77%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +3.68%  "

# Row 12 - Dogecoin
Set-TextCell "D12" "0.0833"
$ws.Range("E12").Value = "  -1.85%  "

# Row 13 - Chainlink
Set-TextCell "D13" "19.99"
$ws.Range("E13").Value = "  +2.81%  "

# Row 14 - Polkadot
Set-TextCell "D14" "7.51"
$ws.Range("E14").Value = "  -1.27%  "

# Row 15 - Wrapped liquid staked Ether 2.0
Set-TextCell "D15" "3.197.52"
$ws.Range("E15").Value = "  -0.64%  "

# Row 16 - Wrapped Ether
Set-TextCell "D16" "2.782.27"
$ws.Range("E16").Value = "  +0.73%  "

# Row 17 - Polygon
Set-TextCell "D17" "0.930"
$ws.Range("E17").Value = "  +0.89%  "

# Row 18 - Wrapped BTC (only the price moved, Volume(1h) is unchanged)
Set-TextCell "D18" "51.131.04"

# Row 19 - Uniswap
$ws.Range("E19").Value = "  +4.03%  "

# Row 20 - ImmutableX
$ws.Range("E20").Value = "  -2.10%  "

# Row 21 - Internet Computer (DFINITY)
$ws.Range("E21").Value = "  -0.13%  "

# Row 22 - Shiba Inu
Set-TextCell "D22" "0.0₃0961"
$ws.Range("E22").Value = "  -1.44%  "

# Row 23 - Litecoin
Set-TextCell "D23" "69.69"
$ws.Range("E23").Value = "  +0.06%  "

# Row 24 - Bitcoin Cash (only the price moved, Volume(1h) is unchanged)
Set-TextCell "D24" "265.29"

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -1.20%  "

# Row 27 - Ethereum Classic
$ws.Range("E27").Value = "  -2.29%  "

# Row 28 - Kaspa
$ws.Range("E28").Value = "  +12.21%  "

# Row 29 - Cosmos
Set-TextCell "D29" "10.16"
$ws.Range("E29").Value = "  +0.36%  "

# Row 30 - Toncoin
$ws.Range("E30").Value = "  +0.77%  "

# Row 31 - OKB
Set-TextCell "D31" "51.78"
$ws.Range("E31").Value = "  +0.72%  "

# Row 32 - Injective Protocol
Set-TextCell "D32" "34.56"
$ws.Range("E32").Value = "  +1.34%  "

# Row 33 - Filecoin
Set-TextCell "D33" "6.03"
$ws.Range("E33").Value = "  +5.58%  "

# Row 34 / 35 - VeChain and RenderToken swapped places in the source ranking
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D34" "5.57"
$ws.Range("E34").Value = "  +2.80%  "

$ws.Range("B35").Value = "VeChain"
$ws.Range("C35").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D35" "0.0443"
$ws.Range("E35").Value = "  -4.33%  "

# Row 36 - Hedera
$ws.Range("E36").Value = "  -0.84%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  +0.05%  "

# Row 38 - Celestia
Set-TextCell "D38" "18.33"
$ws.Range("E38").Value = "  +0.46%  "

# Row 39 - Lido DAO Token
Set-TextCell "D39" "3.13"
$ws.Range("E39").Value = "  -2.34%  "

# Row 40 - ARBITRUM
Set-TextCell "D40" "1.95"
$ws.Range("E40").Value = "  -2.26%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  -1.07%  "

# Row 42 - Stellar
$ws.Range("E42").Value = "  -0.40%  "

# Row 43 - Monero
Set-TextCell "D43" "120.41"
$ws.Range("E43").Value = "  -2.90%  "

# Row 44 - EnergySwap
Set-TextCell "D44" "22.15"
$ws.Range("E44").Value = "  +1.00%  "

# Row 45 - WEMIX Token
$ws.Range("E45").Value = "  -2.59%  "

# Row 46 - Maker
Set-TextCell "D46" "2.087.97"
$ws.Range("E46").Value = "  +1.01%  "

# Row 47 - NEAR Protocol
$ws.Range("E47").Value = "  -0.31%  "

# Row 48 - ApeX Protocol
$ws.Range("E48").Value = "  -1.03%  "

# Row 49 - THORChain
Set-TextCell "D49" "5.49"
$ws.Range("E49").Value = "  -3.58%  "

# Row 50 - SEI
Set-TextCell "D50" "0.918"
$ws.Range("E50").Value = "  -0.78%  "
